$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays formatted as text, matching the source data
# which stores prices as plain strings (not locale-parsed numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = "59.539.67"
$ws.Range("E2").Value2 = "  -0.39%  "

$ws.Range("D3").Value2 = "2.650.34"
$ws.Range("E3").Value2 = "  +0.05%  "

$ws.Range("E4").Value2 = "  -0.23%  "

$ws.Range("D5").Value2 = "518.53"
$ws.Range("E5").Value2 = "  +0.16%  "

$ws.Range("D6").Value2 = "146.52"
$ws.Range("E6").Value2 = "  -0.05%  "

$ws.Range("D7").Value2 = "0.997"
$ws.Range("E7").Value2 = "  +0.17%  "

$ws.Range("D8").Value2 = "0.574"
$ws.Range("E8").Value2 = "  +0.26%  "

$ws.Range("D9").Value2 = "2.658.97"
$ws.Range("E9").Value2 = "  -0.66%  "

$ws.Range("D10").Value2 = "6.32"
$ws.Range("E10").Value2 = "  -2.61%  "

$ws.Range("E11").Value2 = "  -1.31%  "

$ws.Range("E12").Value2 = "  -0.85%  "

$ws.Range("E13").Value2 = "  +0.96%  "

$ws.Range("D14").Value2 = "3.114.31"
$ws.Range("E14").Value2 = "  -0.12%  "

$ws.Range("D15").Value2 = "59.521.02"
$ws.Range("E15").Value2 = "  -0.37%  "

$ws.Range("D16").Value2 = "21.17"
$ws.Range("E16").Value2 = "  -0.66%  "

$ws.Range("D17").Value2 = "0.0000138"
$ws.Range("E17").Value2 = "  -0.64%  "

$ws.Range("D18").Value2 = "2.645.36"
$ws.Range("E18").Value2 = "  -1.71%  "

$ws.Range("D19").Value2 = "350.67"
$ws.Range("E19").Value2 = "  +1.05%  "

$ws.Range("D20").Value2 = "4.52"
$ws.Range("E20").Value2 = "  -2.34%  "

$ws.Range("D21").Value2 = "10.35"
$ws.Range("E21").Value2 = "  -2.05%  "

$ws.Range("D22").Value2 = "6.27"
$ws.Range("E22").Value2 = "  +1.12%  "

$ws.Range("D23").Value2 = "0.999"
$ws.Range("E23").Value2 = "  +0.06%  "

$ws.Range("D24").Value2 = "62.94"
$ws.Range("E24").Value2 = "  +2.85%  "

$ws.Range("E25").Value2 = "  -1.95%  "

$ws.Range("D26").Value2 = "0.166"
$ws.Range("E26").Value2 = "  +2.26%  "

$ws.Range("D27").Value2 = "0.994"

$ws.Range("D28").Value2 = "0.0₃0811"
$ws.Range("E28").Value2 = "  -1.53%  "

$ws.Range("D29").Value2 = "7.18"
$ws.Range("E29").Value2 = "  -0.61%  "

$ws.Range("D30").Value2 = "0.999"
$ws.Range("E30").Value2 = "  +0.09%  "

$ws.Range("E31").Value2 = "  +0.40%  "

$ws.Range("E32").Value2 = "  -0.20%  "

$ws.Range("D33").Value2 = "18.92"
$ws.Range("E33").Value2 = "  -0.88%  "

$ws.Range("D34").Value2 = "150.25"
$ws.Range("E34").Value2 = "  +0.30%  "

$ws.Range("D35").Value2 = "4.08"
$ws.Range("E35").Value2 = "  +0.09%  "

$ws.Range("D36").Value2 = "0.940"
$ws.Range("E36").Value2 = "  -11.90%  "

$ws.Range("D37").Value2 = "1.18"
$ws.Range("E37").Value2 = "  +1.50%  "

$ws.Range("D38").Value2 = "0.871"
$ws.Range("E38").Value2 = "  -0.42%  "

$ws.Range("B39").Value2 = "Stacks"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value2 = "1.49"
$ws.Range("E39").Value2 = "  +3.37%  "

$ws.Range("B40").Value2 = "OKB"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value2 = "36.70"
$ws.Range("E40").Value2 = "  -0.09%  "

$ws.Range("D41").Value2 = "3.69"
$ws.Range("E41").Value2 = "  -1.54%  "

$ws.Range("D42").Value2 = "280.05"
$ws.Range("E42").Value2 = "  -1.10%  "

$ws.Range("D43").Value2 = "0.997"
$ws.Range("E43").Value2 = "  +0.33%  "

$ws.Range("D44").Value2 = "0.0989"
$ws.Range("E44").Value2 = "  -0.58%  "

$ws.Range("D45").Value2 = "19.76"
$ws.Range("E45").Value2 = "  -0.58%  "

$ws.Range("D46").Value2 = "0.603"
$ws.Range("E46").Value2 = "  -3.27%  "

$ws.Range("D47").Value2 = "2.080.73"
$ws.Range("E47").Value2 = "  +4.01%  "

$ws.Range("D48").Value2 = "0.0531"
$ws.Range("E48").Value2 = "  -2.64%  "

$ws.Range("E49").Value2 = "  -0.67%  "

$ws.Range("B50").Value2 = "RenderToken"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value2 = "4.75"
$ws.Range("E50").Value2 = "  +0.18%  "

$ws.Range("B51").Value2 = "WhiteBITCoin"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value2 = "10.30"
$ws.Range("E51").Value2 = "  +0.28%  "
